$wb = $excel.ActiveWorkbook

# --- Sheet: Swing Trade - GTD -> add row 9 ---
$wsGtd = $wb.Worksheets.Item("Swing Trade - GTD")
$wsGtd.Range("A8:E8").Copy()
$wsGtd.Range("A9:E9").PasteSpecial(-4122)
$wsGtd.Range("A9").Value = "b7"
$wsGtd.Range("B9").Value = "Enviar um NewOrderSingle na compra de um ativo da BMF, com validade para amanha"
$wsGtd.Range("C9").Value = "Receber um ExecutionReport com ExecType = PendingNew, e logo após receber um outro com o ExecType = Rejected e a tag Text informando o motivo"
$wsGtd.Rows.Item(9).RowHeight = 30

# --- Sheet: Swing Trade - Stop -> add row 13 ---
$wsStop = $wb.Worksheets.Item("Swing Trade - Stop")
$wsStop.Range("A12:E12").Copy()
$wsStop.Range("A13:E13").PasteSpecial(-4122)
$wsStop.Range("A13").Value = "c11"
$wsStop.Range("B13").Value = "Enviar um NewOrderSingle (OrdType = Stop) na compra de 100 a com o stop em R$2 e sem a tag price price,  com validade para o dia"
$wsStop.Range("C13").Value = "Receber um ExecutionReport com ExecType = PendingNew, e logo após receber um outro ExecutionReport com o ExecType = New e OrdType = Stop. A tag Price devera informara o preco definido para o Stop."
$wsStop.Rows.Item(13).RowHeight = 45

# --- Sheet: DayTrade - Stop -> add row 15 ---
$wsDtStop = $wb.Worksheets.Item("DayTrade - Stop")
$wsDtStop.Range("A14:E14").Copy()
$wsDtStop.Range("A15:E15").PasteSpecial(-4122)
$wsDtStop.Range("A15").Value = "g13"
$wsDtStop.Range("B15").Value = "Enviar um NewOrderSingle (OrdType = Stop) na compra de 100 a com o stop em R$2 e sem a tag price price,  com validade para o dia e TargetStrategy = 1002"
$wsDtStop.Range("C15").Value = "Receber um ExecutionReport com ExecType = PendingNew, e logo após receber um outro ExecutionReport com o ExecType = New e OrdType = Stop. A tag Price devera informara o preco definido para o Stop."
$wsDtStop.Rows.Item(15).RowHeight = 45

# --- Selections / active sheet ---
$wsStop.Range("C13").Select()
$wsDtStop.Range("B15").Select()

$wsGtd.Activate()
$wsGtd.Range("B9").Select()
